$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.796.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.305.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "554.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.295.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.835.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "601.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.858.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.308.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.898"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "566.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.692.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.36%  "
$ws.Range("B41").Value = "CoreDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.128"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0707"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.338"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0417"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("E51").Value = "  +0.08%  "
